# Append 13 new data rows (22-34) to the community smells dataset sheet,
# matching the rows added upstream for the "anago" and "ranking" repos.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: id, repositoryUrl, repositoryName, repositoryAuthor, startingDate,
#             OSE, BCE, PDE, SV, OS, SD, RS, TFS, UI, TC
$newRows = @(
    @(21, "https://github.com/Hironsan/anago", "anago", "Hironsan", "06/26/2017", "0","0","1","1","0","1","0","0","0","0"),
    @(22, "https://github.com/Hironsan/anago", "anago", "Hironsan", "06/26/2017", "0","0","1","1","0","1","0","0","0","0"),
    @(23, "https://github.com/Hironsan/anago", "anago", "Hironsan", "06/26/2017", "0","0","1","1","0","1","0","0","0","0"),
    @(24, "https://github.com/tensorflow/ranking", "ranking", "tensorflow", "12/03/2018", "0","0","0","1","0","1","0","0","1","0"),
    @(25, "https://github.com/tensorflow/ranking", "ranking", "tensorflow", "12/03/2018", "0","0","0","1","0","1","0","0","1","0"),
    @(26, "https://github.com/tensorflow/ranking", "ranking", "tensorflow", "12/03/2018", "0","0","0","1","0","1","0","0","1","0"),
    @(27, "https://github.com/tensorflow/ranking", "ranking", "tensorflow", "12/03/2018", "0","0","0","1","0","1","0","0","1","0"),
    @(28, "https://github.com/tensorflow/ranking", "ranking", "tensorflow", "12/03/2018", "0","0","0","1","1","1","0","0","1","0"),
    @(29, "https://github.com/tensorflow/ranking", "ranking", "tensorflow", "12/03/2018", "0","0","0","1","0","1","0","0","1","0"),
    @(30, "https://github.com/tensorflow/ranking", "ranking", "tensorflow", "12/03/2018", "0","0","0","1","1","1","0","0","1","0"),
    @(31, "https://github.com/tensorflow/ranking", "ranking", "tensorflow", "12/03/2018", "0","0","0","1","0","1","0","0","1","0"),
    @(32, "https://github.com/tensorflow/ranking", "ranking", "tensorflow", "12/03/2018", "0","0","0","1","1","1","0","0","1","0"),
    @(33, "https://github.com/tensorflow/ranking", "ranking", "tensorflow", "12/03/2018", "0","0","0","1","0","1","0","0","1","0")
)

$startRow = 22
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $entry = $newRows[$i]

    # Column A: numeric id.
    $ws.Cells.Item($r, 1).Value = $entry[0]

    # Columns B-O: text values. The 0/1 flag columns look numeric, so force
    # them to stay text (matching the source data) by applying a text
    # number format before the assignment, then clearing the format again
    # so the cell keeps the default style like its neighbours.
    for ($j = 1; $j -lt $entry.Count; $j++) {
        $col = $j + 1
        $cell = $ws.Cells.Item($r, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $entry[$j]
        $cell.ClearFormats()
    }
}
